# Round the emmean/SE/lower.CL/upper.CL columns (B, C, E, F) for the data
# rows (2 through 28) to 3 decimal places, leaving the age_group (A),
# df (D) and measure (G) columns, as well as the header row, untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 28; $row++) {
    foreach ($col in 2, 3, 5, 6) {
        $cell = $ws.Cells.Item($row, $col)
        $cell.Value2 = [Math]::Round($cell.Value2, 3)
    }
}
